$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content -----------------------------------------------------
# Order matters: the shared-strings table is built in the order cells are
# written, and we want to reproduce the author's apparent edit sequence
# (fix the two existing dates, type the new Day-3 content, add the new
# Day 4-6 placeholder rows, then reword the existing Day1/Day2 summaries).
$ws.Range("A1").Value = "Day 1 (26/05/2019)"
$ws.Range("A2").Value = "Day 2 (27/05/2019)"
$ws.Range("B3").Value = "Learning how to install scss with node`nLearning one CSS property => clip-path `n"
$ws.Range("A3").Value = "Day 3 (28/05/2019)"
$ws.Range("A4").Value = "Day 4 (29/05/2019)"
$ws.Range("A5").Value = "Day 5 (30/05/2019)"
$ws.Range("A6").Value = "Day 6 (31/05/2019)"
$ws.Range("B2").Value = "Learning using Sass Mixins Extends and Functions"
$ws.Range("B1").Value = "Learning using Sass Variables and Nesting"

# --- Cell formatting ----------------------------------------------------
# Each distinct alignment combination is finished on a single cell first
# (so no partially-applied intermediate style leaks into the style table),
# then copied (format only) onto the remaining cells that share it.

# B3: left/top, wrap text
$b3 = $ws.Range("B3")
$b3.HorizontalAlignment = -4131   # xlLeft
$b3.VerticalAlignment = -4160     # xlTop
$b3.WrapText = $true

# Column A (A1:A6): centered both ways
$a1 = $ws.Range("A1")
$a1.HorizontalAlignment = -4108   # xlCenter
$a1.VerticalAlignment = -4108     # xlCenter
$a1.Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)   # xlPasteFormats

# B1 / B2: left/top, no wrap
$b2 = $ws.Range("B2")
$b2.HorizontalAlignment = -4131   # xlLeft
$b2.VerticalAlignment = -4160     # xlTop
$b2.Copy()
$ws.Range("B1").PasteSpecial(-4122)      # xlPasteFormats

$ws.Application.CutCopyMode = $false

# Row 3 needs extra height for the two-line Day-3 note.
$ws.Rows(3).RowHeight = 45

# --- Column widths --------------------------------------------------
$ws.Columns("A").ColumnWidth = 24.666666666666668
$ws.Columns("B").ColumnWidth = 71

# --- Page setup -------------------------------------------------------
$ws.PageSetup.Orientation = 1    # xlPortrait

# --- Selection ----------------------------------------------------------
$ws.Range("B9").Select() | Out-Null
